$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update existing GDP per Capita values for years 1952-2010 (rows 2-60)
$gdpValues = @(
    "3934",
    "4106",
    "4623",
    "5002",
    "4889",
    "5432",
    "6005",
    "6304",
    "7165",
    "7760",
    "8027",
    "8871",
    "9624",
    "9589",
    "10029",
    "10264",
    "10822",
    "11775",
    "12803",
    "13657",
    "14269",
    "15079",
    "16420",
    "17205",
    "17462",
    "18595",
    "20014",
    "21367",
    "21151",
    "20929",
    "20535",
    "20451",
    "20650",
    "20663",
    "21122",
    "20851",
    "20083",
    "19837",
    "18093",
    "16401.95914321",
    "15462.7101618003",
    "15872.1202626296",
    "16664.70579348",
    "17816.8010076464",
    "18376.2462010584",
    "19258.9458702082",
    "19851.0571234378",
    "20793.9180282206",
    "21501.3367278983",
    "22006.9537877351",
    "22725.3119846003",
    "23257.421804062",
    "24150.5840374106",
    "24966.9387886715",
    "26182.7796675607",
    "27832.7637848008",
    "28473.5957291968",
    "25905.5126735662",
    "26000.9206515415"
)

for ($i = 0; $i -lt $gdpValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Formula = "'" + $gdpValues[$i]
}

# Append new rows for years 2011-2016 (rows 61-66)
$newYears = @(2011,2012,2013,2014,2015,2016)
$newValues = @("26004","25252","24944","25693","26268","26908")

for ($i = 0; $i -lt $newYears.Length; $i++) {
    $row = 61 + $i
    $ws.Cells.Item($row, 1).Value = 705
    $ws.Cells.Item($row, 2).Value = "Slovenia"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $newYears[$i]
    $ws.Cells.Item($row, 5).Formula = "'" + $newValues[$i]
}
